$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: new "D/M/YYYY" format sample row
$ws.Range("A7").Value = "D/M/YYYY"
$ws.Range("B7").Value = 32
$ws.Range("C7").Value = 60
$ws.Range("D7").Value = 44228
$ws.Range("B7:D7").NumberFormat = "d/m/yyyy;@"

# Row 8: d/mm/yy
$ws.Range("B8").Value = 32
$ws.Range("B8").NumberFormat = "d/mm/yy;@"

# Row 9: d/m/yy
$ws.Range("B9").Value = 32
$ws.Range("B9").NumberFormat = "d/m/yy;@"

# Row 10: dd/mm/yy
$ws.Range("B10").Value = 32
$ws.Range("B10").NumberFormat = "dd/mm/yy;@"

# Row 11: [$-C09]dd-mmm-yy
$ws.Range("B11").Value = 32
$ws.Range("B11").NumberFormat = "[$-C09]dd-mmm-yy;@"

# Row 12: [$-C09]dd-mmmm-yyyy
$ws.Range("B12").Value = 32
$ws.Range("B12").NumberFormat = "[$-C09]dd-mmmm-yyyy;@"

# Row 13: [$-C09]d mmmm yyyy
$ws.Range("B13").Value = 32
$ws.Range("B13").NumberFormat = "[$-C09]d mmmm yyyy;@"

# Update selection to follow the newly added data (now on row 14)
$ws.Range("B14").Select()

# Touch page setup so dpi settings get (re)written
$ws.PageSetup.PrintQuality = 300
